# PM08 Tidsregistrering for Benjamin.xlsx
# Adds three new logged time-tracking entries (rows 27-29, pushing the
# old placeholder rows down), extends the running-total formulas in
# columns G/H to the newly created row 39, restores the blank formatted
# cell in F31, appends one more blank trailing row (54), and leaves the
# active selection on C18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 27: new task entry -------------------------------------------------
$ws.Cells.Item(27, 1).Value = "Lav OC0803 angivSaldoafskrivning"   # A27 (Opgavebeskrivelse)
$ws.Cells.Item(27, 2).Value = "System Analyst "                    # B27 (Rolle)
$ws.Cells.Item(27, 3).Value = 43893                                 # C27 (Dato)
$ws.Cells.Item(27, 4).Value = 0.53472222222222221                  # D27 (Starttid)
$ws.Cells.Item(27, 5).Value = 0.60416666666666663                  # E27 (Sluttid)
$ws.Cells.Item(27, 6).Value = "1t 30min"                            # F27 (Estimeret tidsforbrug)

# --- Row 28: new task entry -------------------------------------------------
$ws.Cells.Item(28, 1).Value = "Review OC0802 og DCD0802"
$ws.Cells.Item(28, 2).Value = "Reviewer"
$ws.Cells.Item(28, 3).Value = 43894
$ws.Cells.Item(28, 4).Value = 0.60416666666666663
$ws.Cells.Item(28, 5).Value = 0.625
# H28's running-total formula is removed for this row in the final workbook.
$ws.Cells.Item(28, 8).ClearContents()

# --- Row 29: new task entry -------------------------------------------------
$ws.Cells.Item(29, 1).Value = "Lav SD0802 og DCD0802 angivLineaerAfskrivning"
$ws.Cells.Item(29, 2).Value = "System Analyst "
$ws.Cells.Item(29, 3).Value = 43895
$ws.Cells.Item(29, 4).Value = 0.625
$ws.Cells.Item(29, 5).Value = 0.67013888888888884

# --- Row 31 regains its blank, formatted F cell -----------------------------
$ws.Range("F30").Copy() | Out-Null
$ws.Range("F31").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# --- New row 39: extend the diff-formula (G) and running-total (H) ----------
$ws.Range("G39").Formula = "=E39-D39"
$ws.Range("H39").Formula = "=SUM(G`$5:G39)"
$ws.Range("G38").Copy() | Out-Null
$ws.Range("G39").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H38").Copy() | Out-Null
$ws.Range("H39").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# --- New trailing blank row 54 (matches formatting of row 53) --------------
$ws.Range("C53").Copy() | Out-Null
$ws.Range("C54").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# --- Final selection ---------------------------------------------------------
$ws.Range("C18").Select()
